$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sd_marca")
$ws.Range("B2").Value = "LORETO"
$v = $ws.Range("B2").Value2
Write-Host ("New Value2 is: " + $v)
